$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.976.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.424.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.59%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.008.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.434.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.008.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '391.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.04%  '
$ws.Range("E22").Value = '  -1.34%  '
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("E26").Value = '  +2.72%  '
$ws.Range("E27").Value = '  +3.73%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("E31").Value = '  +3.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.57'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.30'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.97'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '167.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.456.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '28.33'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0756'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("E41").Value = '  +1.28%  '
$ws.Range("E42").Value = '  +1.88%  '
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("E44").Value = '  +4.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.517.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = '  -0.41%  '
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.808'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.79%  '
